# CDIO3 Chancekort - implement a few Chancekort card-text updates
# (commit: "Nogle Chancekort funktioner er implementeret")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ChanceKort")

# Bilen (car) card: clarify that the card is handed over before re-drawing
$ws.Range("C2").Value = "Dette chance kort er givet til Bilen. Tag et chancekort mere. Bil: På din næste tur skal du drøne frem til et hvilket som helst ledigt felt og købe det. Hvis det ikke er nogen ledige felter skal du købe et fra en anden spiller!"

# "Ryk 1 felt frem eller tag et chancekort mere" row: the choice answer was "?" - clarified to "Nej"
$ws.Range("B6").Value = "Nej"

# Skibet (ship) card: clarify that the card is handed over before re-drawing
$ws.Range("C7").Value = "Dette kort er givet til skibet.Tag et chancekort mere. Skib: På den næste skal du sejle frem til hvilket som helst ledigt felt og købe det. Hvis der ikke er nogen ledige felter, skal du købe et af en anden spiller."

# Katten (cat) card: clarify that the card is handed over before re-drawing
$ws.Range("C13").Value = "Dette kort er givet til Katten . Tag et chancekort mere. Kat: På den næste skal du sejle frem til hvilket som helst ledigt felt og købe det. Hvis der ikke er nogen ledige felter, skal du købe et af en anden spiller."

# Hunden (dog) card: clarify that the card is handed over before re-drawing
$ws.Range("C14").Value = "Dette kort er givet til Hunden. Tag et chancekort mere. Hund: På den næste skal du sejle frem til hvilket som helst ledigt felt og købe det. Hvis der ikke er nogen ledige felter, skal du købe et af en anden spiller."
